{"js": "// Update the worksheet date header and the 25 multiplication problems\n// (5 populated rows x 5 columns) inside the single table, in place,\n// preserving all existing run/paragraph formatting.\n\n// 1) Title paragraph: \"2024-05-23 Thursday\" -> \"2024-05-24 Friday\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].getRange().insertText(\"2024-05-24 Friday\", \"Replace\");\n\n// 2) Table values: only 5 of the 20 rows actually contain problems\n// (rows 0, 4, 9, 14, 19); every row has 5 columns. Values are addressed\n// by stable (row, col) position rather than by searching for the old\n// text, because several old/new strings collide with each other\n// (e.g. \"552\u00d78=\" becomes \"345\u00d76=\", while a *different* cell's original\n// \"345\u00d76=\" becomes \"969\u00d72=\") which would make a plain text search\n// ambiguous after earlier replacements land.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"114\u00d79=\", \"345\u00d76=\", \"725\u00d77=\", \"961\u00d73=\", \"205\u00d78=\"],\n  4: [\"824\u00d72=\", \"217\u00d72=\", \"343\u00d78=\", \"810\u00d75=\", \"640\u00d77=\"],\n  9: [\"285\u00d74=\", \"188\u00d75=\", \"412\u00d73=\", \"986\u00d78=\", \"683\u00d75=\"],\n  14: [\"682\u00d76=\", \"361\u00d72=\", \"543\u00d72=\", \"820\u00d79=\", \"691\u00d79=\"],\n  19: [\"822\u00d74=\", \"838\u00d73=\", \"163\u00d73=\", \"969\u00d72=\", \"991\u00d76=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 multiplication problems\n# (5 populated rows x 5 columns) inside the single table, in place,\n# preserving all existing run/paragraph formatting.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2024-05-23 Thursday\" -> \"2024-05-24 Friday\"\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-24 Friday\"\n\n# 2) Table values: only 5 of the 20 rows actually contain problems\n# (1-based rows 1, 5, 10, 15, 20); every row has 5 columns. Values are\n# addressed by stable (row, col) position rather than by Find/Replace\n# on the old text, because several old/new strings collide with each\n# other (e.g. \"552\u00d78=\" becomes \"345\u00d76=\", while a *different* cell's\n# original \"345\u00d76=\" becomes \"969\u00d72=\"), which would make plain text\n# search-and-replace ambiguous once earlier replacements have landed.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"114\u00d79=\", \"345\u00d76=\", \"725\u00d77=\", \"961\u00d73=\", \"205\u00d78=\")\n    5  = @(\"824\u00d72=\", \"217\u00d72=\", \"343\u00d78=\", \"810\u00d75=\", \"640\u00d77=\")\n    10 = @(\"285\u00d74=\", \"188\u00d75=\", \"412\u00d73=\", \"986\u00d78=\", \"683\u00d75=\")\n    15 = @(\"682\u00d76=\", \"361\u00d72=\", \"543\u00d72=\", \"820\u00d79=\", \"691\u00d79=\")\n    20 = @(\"822\u00d74=\", \"838\u00d73=\", \"163\u00d73=\", \"969\u00d72=\", \"991\u00d76=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
